# SCD0018-013 - Admin WEM mengajukan data Non Sales.xlsx
# Update TC_ID Excel SCD0017 until SCD0025 and Update TC_ID Solution SCD0006 until SCD0025
#
# Changes applied:
#  1. Rename worksheet "SCD0290" -> "SCD0018"
#  2. Update TC_ID cells (column B, rows 2 and 3) from "DGS-305" to "SCD0018-013"
#  3. Drop the wrapText formatting on the TC_ID (B) and TEST_SCENARIO_DESC (C) cells for
#     rows 2 & 3 (they move to the already-used "no-wrap" style)
#  4. Re-fit column B so the new, longer TC_ID value is not clipped
#  5. Move the active selection from E4 to B4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet
$ws.Name = "SCD0018"

# 2 & 3. Update the TC_ID values and drop wrap-text formatting for B2:C3
$ws.Range("B2").Value = "SCD0018-013"
$ws.Range("B2").WrapText = $false
$ws.Range("B2").Font.Size = 10

$ws.Range("C2").WrapText = $false
$ws.Range("C2").Font.Size = 10

$ws.Range("B3").Value = "SCD0018-013"
$ws.Range("B3").WrapText = $false
$ws.Range("B3").Font.Size = 10

$ws.Range("C3").WrapText = $false
$ws.Range("C3").Font.Size = 10

# 4. Re-fit column B to the new (wider) TC_ID text
$ws.Columns("B:B").AutoFit()

# 5. Update the active cell/selection
[void]$ws.Range("B4").Select()
